$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.768.48'
$ws.Range('E2').Value = '  +0.80%  '
$ws.Range('D3').Value = '2.491.25'
$ws.Range('E3').Value = '  +0.04%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '''586.54'
$ws.Range('E5').Value = '  +0.34%  '
$ws.Range('E6').Value = '  +2.65%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  +0.17%  '
$ws.Range('E9').Value = '  +3.28%  '
$ws.Range('E10').Value = '  -0.19%  '
$ws.Range('E12').Value = '  +0.06%  '
$ws.Range('D13').Value = '2.946.55'
$ws.Range('E13').Value = '  +0.92%  '
$ws.Range('D14').Value = '''25.62'
$ws.Range('E14').Value = '  +0.60%  '
$ws.Range('D15').Value = '67.678.13'
$ws.Range('E15').Value = '  +0.91%  '
$ws.Range('E16').Value = '  +0.49%  '
$ws.Range('D17').Value = '2.497.06'
$ws.Range('E17').Value = '  +0.78%  '
$ws.Range('E18').Value = '  -0.87%  '
$ws.Range('D19').Value = '''7.47'
$ws.Range('D20').Value = '''351.50'
$ws.Range('E20').Value = '  +0.28%  '
$ws.Range('E21').Value = '  +2.40%  '
$ws.Range('E22').Value = '  +0.11%  '
$ws.Range('E23').Value = '  +3.07%  '
$ws.Range('E24').Value = '  +1.64%  '
$ws.Range('E25').Value = '  -3.30%  '
$ws.Range('E26').Value = '  -1.60%  '
$ws.Range('D27').Value = '2.619.87'
$ws.Range('D28').Value = '''0.992'
$ws.Range('E28').Value = '  -0.71%  '
$ws.Range('D29').Value = '0.0₃0903'
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('D30').Value = '''503.61'
$ws.Range('E30').Value = '  -1.19%  '
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('E32').Value = '  +1.87%  '
$ws.Range('E33').Value = '  +0.17%  '
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('E35').Value = '  +3.05%  '
$ws.Range('D36').Value = '''163.60'
$ws.Range('E36').Value = '  +2.21%  '
$ws.Range('D37').Value = '''18.65'
$ws.Range('E37').Value = '  -0.30%  '
$ws.Range('D38').Value = '''18.32'
$ws.Range('E38').Value = '  +0.37%  '
$ws.Range('E39').Value = '  -0.13%  '
$ws.Range('E41').Value = '  +2.85%  '
$ws.Range('E42').Value = '  +0.13%  '
$ws.Range('D43').Value = '''4.85'
$ws.Range('E43').Value = '  +0.33%  '
$ws.Range('E44').Value = '  +2.36%  '
$ws.Range('D45').Value = '''144.47'
$ws.Range('E45').Value = '  +1.14%  '
$ws.Range('E46').Value = '  +1.94%  '
$ws.Range('D47').Value = '''0.514'
$ws.Range('E47').Value = '  -0.06%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₆0254'
$ws.Range('E48').Value = '  +1.29%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '''0.0742'
$ws.Range('E49').Value = '  +1.67%  '
$ws.Range('B50').Value = 'Optimism'
$ws.Range('C50').Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range('D50').Value = '''1.58'
$ws.Range('E50').Value = '  +0.86%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').Value = '''0.585'
$ws.Range('E51').Value = '  +0.18%  '
